$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# B2 holds the "Cases" Cypher query used to build the CasesTab export. It
# used to return a trailing `Cohort` column (coalesce(co.cohort_description,
# '') AS `Cohort`); that final RETURN column is being dropped, so update the
# cell text to the same query minus that last line.
$newQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n`nMATCH (c)<--(diag:diagnosis)`n MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)`n`tWHERE s.clinical_study_designation IN ['UBC01'] and diag.stage_of_disease in ['T3N1M0'] OPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newQuery

# The row height is driven by the wrapped query text; with one less line the
# row is a bit shorter than before.
$ws.Rows.Item(2).RowHeight = 290

# Restore the view: scrolled to the top of the sheet with B2 selected.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B2").Select()
